$wb = $excel.ActiveWorkbook

$upper = $wb.Worksheets.Item("Upper")
$lower = $wb.Worksheets.Item("Lower")
$mtul  = $wb.Worksheets.Item("Mid to upper lower")

$rows = @(5,6,7,9,10,11,12,13)

foreach ($r in $rows) {
    $upper.Range("B$r").Formula = "=Mid!B$r+'Mid to upper lower'!B$r"
    $lower.Range("B$r").Formula = "=Mid!B$r-'Mid to upper lower'!B$r"
}

# "Mid to upper lower": intake = 0 always has RR = 0 (was blank before)
$zeroRows = @(5,6,7,9,10,11,13)
foreach ($r in $zeroRows) {
    $c = $mtul.Range("B$r")
    $c.Value = 0
    $c.Font.Color = 0
}
$mtul.Range("B12").Value = 0

# Selections / active sheet, to mirror the recorded view state.
$mtul.Range("B17").Select()
$upper.Range("C19").Select()

$lower.Activate()
$lower.Range("B4:B13").Select()
